$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Correct the product name / short-name code value (typo fix: add missing dash)
$newValue = "293-MS-EPP-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$ws1.Range("B1").Value = $newValue
$ws2.Range("B1").Value = $newValue

# Move the selection on the input sheet to B1
[void]$ws1.Range("B1").Select()

# Activate the output sheet and select its B1 cell, making it the active tab
[void]$ws2.Activate()
[void]$ws2.Range("B1").Select()
